$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.528.75"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.515.25"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.513.46"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.55%  "
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.974.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.387.21"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.82"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.515.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.25"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.645.49"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0887"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "461.92"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.24"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.45"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.318"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -7.39%  "
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.92"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.519"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.46"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0732"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  -3.28%  "
